$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28:100 down to 29:101.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44497
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112003
$ws.Range("G28").Value = "Ajo"
$ws.Range("H28").Value = "Chino"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 270
$ws.Range("K28").Value = 17000
$ws.Range("L28").Value = 18000
$ws.Range("M28").Value = 17444
$ws.Range("N28").Value = "`$/caja 10 kilos"
$ws.Range("O28").Value = "China"
$ws.Range("P28").Value = 1744
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = "Hortaliza"
